$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.038172232093932
$ws.Cells.Item(2, 4).Value = 1.039141639893988
$ws.Cells.Item(2, 5).Value = 1.0514363096429
$ws.Cells.Item(2, 6).Value = 1.058297005690545
$ws.Cells.Item(2, 9).Value = 1.032951788875909
$ws.Cells.Item(2, 10).Value = 1.043271113561038
$ws.Cells.Item(2, 11).Value = 1.041927670689597
$ws.Cells.Item(2, 12).Value = 1.054187816701905
$ws.Cells.Item(2, 13).Value = 1.061029625515856
$ws.Cells.Item(2, 14).Value = 1.044752678207164

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.039356689080995
$ws.Cells.Item(3, 4).Value = 1.040012921798461
$ws.Cells.Item(3, 5).Value = 1.05264503947401
$ws.Cells.Item(3, 6).Value = 1.059650120503179
$ws.Cells.Item(3, 9).Value = 1.033144031776881
$ws.Cells.Item(3, 10).Value = 1.044099074951224
$ws.Cells.Item(3, 11).Value = 1.042608949108766
$ws.Cells.Item(3, 12).Value = 1.055208173996032
$ws.Cells.Item(3, 13).Value = 1.062195378954268
$ws.Cells.Item(3, 14).Value = 1.045581815397493

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.040122886448336
$ws.Cells.Item(4, 4).Value = 1.040576192680133
$ws.Cells.Item(4, 5).Value = 1.053427350078393
$ws.Cells.Item(4, 6).Value = 1.060526125145609
$ws.Cells.Item(4, 9).Value = 1.033266847303909
$ws.Cells.Item(4, 10).Value = 1.044634082152327
$ws.Cells.Item(4, 11).Value = 1.04304863947593
$ws.Cells.Item(4, 12).Value = 1.055868027236708
$ws.Cells.Item(4, 13).Value = 1.062949608084786
$ws.Cells.Item(4, 14).Value = 1.046117582370187

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.040444943716996
$ws.Cells.Item(5, 4).Value = 1.040812870924597
$ws.Cells.Item(5, 5).Value = 1.053756278501493
$ws.Cells.Item(5, 6).Value = 1.060894506991295
$ws.Cells.Item(5, 9).Value = 1.033318101171732
$ws.Cells.Item(5, 10).Value = 1.044858823255202
$ws.Cells.Item(5, 11).Value = 1.043233212437441
$ws.Cells.Item(5, 12).Value = 1.05614533846819
$ws.Cells.Item(5, 13).Value = 1.063266665853721
$ws.Cells.Item(5, 14).Value = 1.046342642631207

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.040499015548782
$ws.Cells.Item(6, 4).Value = 1.040852603171598
$ws.Cells.Item(6, 5).Value = 1.053811509718306
$ws.Cells.Item(6, 6).Value = 1.06095636640398
$ws.Cells.Item(6, 9).Value = 1.03332668477173
$ws.Cells.Item(6, 10).Value = 1.044896547989099
$ws.Cells.Item(6, 11).Value = 1.043264187066235
$ws.Cells.Item(6, 12).Value = 1.056191894943716
$ws.Cells.Item(6, 13).Value = 1.063319900115568
$ws.Cells.Item(6, 14).Value = 1.046380420938553

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.040127189995709
$ws.Cells.Item(7, 4).Value = 1.040579355660623
$ws.Cells.Item(7, 5).Value = 1.053431745057122
$ws.Cells.Item(7, 6).Value = 1.060531047049088
$ws.Cells.Item(7, 9).Value = 1.03326753364488
$ws.Cells.Item(7, 10).Value = 1.044637085843921
$ws.Cells.Item(7, 11).Value = 1.043051106820323
$ws.Cells.Item(7, 12).Value = 1.055871733039092
$ws.Cells.Item(7, 13).Value = 1.062953844703302
$ws.Cells.Item(7, 14).Value = 1.046120590327368

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.038572572703999
$ws.Cells.Item(8, 4).Value = 1.039436198918873
$ws.Cells.Item(8, 5).Value = 1.05184476890024
$ws.Cells.Item(8, 6).Value = 1.058754205456681
$ws.Cells.Item(8, 9).Value = 1.033017085168613
$ws.Cells.Item(8, 10).Value = 1.043551080667344
$ws.Cells.Item(8, 11).Value = 1.042158148813518
$ws.Cells.Item(8, 12).Value = 1.054532732012383
$ws.Cells.Item(8, 13).Value = 1.061423617274937
$ws.Cells.Item(8, 14).Value = 1.045033042898874

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.03583132675566
$ws.Cells.Item(9, 4).Value = 1.037417906858443
$ws.Cells.Item(9, 5).Value = 1.04904963392638
$ws.Cells.Item(9, 6).Value = 1.055626530191347
$ws.Cells.Item(9, 9).Value = 1.032563668283893
$ws.Cells.Item(9, 10).Value = 1.04163169266788
$ws.Cells.Item(9, 11).Value = 1.040575866168762
$ws.Cells.Item(9, 12).Value = 1.052170221617907
$ws.Cells.Item(9, 13).Value = 1.058726385429146
$ws.Cells.Item(9, 14).Value = 1.043110929148293

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.034002490896422
$ws.Cells.Item(10, 4).Value = 1.036069719582959
$ws.Cells.Item(10, 5).Value = 1.047186996294168
$ws.Cells.Item(10, 6).Value = 1.053543533484919
$ws.Cells.Item(10, 9).Value = 1.032253249778373
$ws.Cells.Item(10, 10).Value = 1.040348195080922
$ws.Cells.Item(10, 11).Value = 1.039515069075992
$ws.Cells.Item(10, 12).Value = 1.050593095502187
$ws.Cells.Item(10, 13).Value = 1.056927600227164
$ws.Cells.Item(10, 14).Value = 1.041825608847543

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033210239751654
$ws.Cells.Item(11, 4).Value = 1.035485300259209
$ws.Cells.Item(11, 5).Value = 1.046380616574343
$ws.Cells.Item(11, 6).Value = 1.052642046114181
$ws.Cells.Item(11, 9).Value = 1.032116901914395
$ws.Cells.Item(11, 10).Value = 1.039791484316384
$ws.Cells.Item(11, 11).Value = 1.039054312941501
$ws.Cells.Item(11, 12).Value = 1.049909659283859
$ws.Cells.Item(11, 13).Value = 1.056148531762935
$ws.Cells.Item(11, 14).Value = 1.041268107489864

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032915906876779
$ws.Cells.Item(12, 4).Value = 1.035268123094971
$ws.Cells.Item(12, 5).Value = 1.046081112365282
$ws.Cells.Item(12, 6).Value = 1.052307260620145
$ws.Cells.Item(12, 9).Value = 1.032065965324646
$ws.Cells.Item(12, 10).Value = 1.039584553296389
$ws.Cells.Item(12, 11).Value = 1.038882952705238
$ws.Cells.Item(12, 12).Value = 1.049655719119038
$ws.Cells.Item(12, 13).Value = 1.055859121361717
$ws.Cells.Item(12, 14).Value = 1.041060882604085

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032979044829883
$ws.Cells.Item(13, 4).Value = 1.035314712773613
$ws.Cells.Item(13, 5).Value = 1.046145356112973
$ws.Cells.Item(13, 6).Value = 1.052379070253333
$ws.Cells.Item(13, 9).Value = 1.032076904565256
$ws.Cells.Item(13, 10).Value = 1.03962894723572
$ws.Cells.Item(13, 11).Value = 1.03891971979734
$ws.Cells.Item(13, 12).Value = 1.049710193865546
$ws.Cells.Item(13, 13).Value = 1.055921202254761
$ws.Cells.Item(13, 14).Value = 1.041105339587902

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033185911251382
$ws.Cells.Item(14, 4).Value = 1.035467350327232
$ws.Cells.Item(14, 5).Value = 1.04635585904434
$ws.Cells.Item(14, 6).Value = 1.052614371288266
$ws.Cells.Item(14, 9).Value = 1.032112697416291
$ws.Cells.Item(14, 10).Value = 1.039774382276727
$ws.Cells.Item(14, 11).Value = 1.0390401526402
$ws.Cells.Item(14, 12).Value = 1.049888670184485
$ws.Cells.Item(14, 13).Value = 1.056124609627509
$ws.Cells.Item(14, 14).Value = 1.041250981163348

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033313361049629
$ws.Cells.Item(15, 4).Value = 1.035561382375828
$ws.Cells.Item(15, 5).Value = 1.046485559583478
$ws.Cells.Item(15, 6).Value = 1.052759356827786
$ws.Cells.Item(15, 9).Value = 1.032134712015002
$ws.Cells.Item(15, 10).Value = 1.039863970499799
$ws.Cells.Item(15, 11).Value = 1.039114326785464
$ws.Cells.Item(15, 12).Value = 1.049998624476484
$ws.Cells.Item(15, 13).Value = 1.056249931622064
$ws.Cells.Item(15, 14).Value = 1.041340696611973

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.034055062344711
$ws.Cells.Item(16, 4).Value = 1.036108491899721
$ws.Cells.Item(16, 5).Value = 1.04724051613878
$ws.Cells.Item(16, 6).Value = 1.053603371781657
$ws.Cells.Item(16, 9).Value = 1.032262257954203
$ws.Cells.Item(16, 10).Value = 1.040385122047153
$ws.Cells.Item(16, 11).Value = 1.039545617859992
$ws.Cells.Item(16, 12).Value = 1.05063844158322
$ws.Cells.Item(16, 13).Value = 1.056979300396066
$ws.Cells.Item(16, 14).Value = 1.041862588254302

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034520215682433
$ws.Cells.Item(17, 4).Value = 1.03645150587532
$ws.Cells.Item(17, 5).Value = 1.047714120357835
$ws.Cells.Item(17, 6).Value = 1.054132922407228
$ws.Cells.Item(17, 9).Value = 1.03234174601395
$ws.Cells.Item(17, 10).Value = 1.040711771856554
$ws.Cells.Item(17, 11).Value = 1.039815773441091
$ws.Cells.Item(17, 12).Value = 1.05103963868386
$ws.Cells.Item(17, 13).Value = 1.057436763682286
$ws.Cells.Item(17, 14).Value = 1.042189701943899

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034791497848945
$ws.Cells.Item(18, 4).Value = 1.036651517917408
$ws.Cells.Item(18, 5).Value = 1.047990380753475
$ws.Cells.Item(18, 6).Value = 1.054441845272674
$ws.Cells.Item(18, 9).Value = 1.032387923412092
$ws.Cells.Item(18, 10).Value = 1.040902209810492
$ws.Cells.Item(18, 11).Value = 1.039973213291894
$ws.Cells.Item(18, 12).Value = 1.051273599135533
$ws.Cells.Item(18, 13).Value = 1.057703576812635
$ws.Cells.Item(18, 14).Value = 1.042380410341575

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034883992380289
$ws.Cells.Item(19, 4).Value = 1.036719706328866
$ws.Cells.Item(19, 5).Value = 1.048084581028237
$ws.Cells.Item(19, 6).Value = 1.054547187689461
$ws.Cells.Item(19, 9).Value = 1.032403637061712
$ws.Cells.Item(19, 10).Value = 1.040967128779765
$ws.Cells.Item(19, 11).Value = 1.040026872947418
$ws.Cells.Item(19, 12).Value = 1.051353365003932
$ws.Cells.Item(19, 13).Value = 1.057794550300234
$ws.Cells.Item(19, 14).Value = 1.04244542150323

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.034470312657507
$ws.Cells.Item(20, 4).Value = 1.036414710163289
$ws.Cells.Item(20, 5).Value = 1.047663305565208
$ws.Cells.Item(20, 6).Value = 1.054076102008296
$ws.Cells.Item(20, 9).Value = 1.032333237004964
$ws.Cells.Item(20, 10).Value = 1.040676734883863
$ws.Cells.Item(20, 11).Value = 1.039786802516184
$ws.Cells.Item(20, 12).Value = 1.05099659932458
$ws.Cells.Item(20, 13).Value = 1.057387684024959
$ws.Cells.Item(20, 14).Value = 1.042154615214688

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033124995798867
$ws.Cells.Item(21, 4).Value = 1.035422405078887
$ws.Cells.Item(21, 5).Value = 1.046293870611613
$ws.Cells.Item(21, 6).Value = 1.052545079168485
$ws.Cells.Item(21, 9).Value = 1.03210216534683
$ws.Cells.Item(21, 10).Value = 1.039731559264803
$ws.Cells.Item(21, 11).Value = 1.039004694109469
$ws.Cells.Item(21, 12).Value = 1.049836115622424
$ws.Cells.Item(21, 13).Value = 1.056064712063612
$ws.Cells.Item(21, 14).Value = 1.041208097337836

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032278819219059
$ws.Cells.Item(22, 4).Value = 1.034797937184593
$ws.Cells.Item(22, 5).Value = 1.045432973663008
$ws.Cells.Item(22, 6).Value = 1.051582850981502
$ws.Cells.Item(22, 9).Value = 1.031955198102879
$ws.Cells.Item(22, 10).Value = 1.039136456818317
$ws.Cells.Item(22, 11).Value = 1.038511707732643
$ws.Cells.Item(22, 12).Value = 1.049106000862192
$ws.Cells.Item(22, 13).Value = 1.055232734212015
$ws.Cells.Item(22, 14).Value = 1.040612149777617

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032727424827148
$ws.Cells.Item(23, 4).Value = 1.035129033396213
$ws.Cells.Item(23, 5).Value = 1.04588934064752
$ws.Cells.Item(23, 6).Value = 1.052092910587068
$ws.Cells.Item(23, 9).Value = 1.032033267869812
$ws.Cells.Item(23, 10).Value = 1.039452011332285
$ws.Cells.Item(23, 11).Value = 1.038773167354795
$ws.Cells.Item(23, 12).Value = 1.049493093932578
$ws.Cells.Item(23, 13).Value = 1.055673798552627
$ws.Cells.Item(23, 14).Value = 1.040928152415189

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034492861791193
$ws.Cells.Item(24, 4).Value = 1.036431336754037
$ws.Cells.Item(24, 5).Value = 1.047686266536193
$ws.Cells.Item(24, 6).Value = 1.054101776561913
$ws.Cells.Item(24, 9).Value = 1.0323370824366
$ws.Cells.Item(24, 10).Value = 1.040692566865987
$ws.Cells.Item(24, 11).Value = 1.039799893654511
$ws.Cells.Item(24, 12).Value = 1.051016047116144
$ws.Cells.Item(24, 13).Value = 1.057409861063485
$ws.Cells.Item(24, 14).Value = 1.042170469680043

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.036540231621706
$ws.Cells.Item(25, 4).Value = 1.037940149884478
$ws.Cells.Item(25, 5).Value = 1.04977209736391
$ws.Cells.Item(25, 6).Value = 1.056434725986676
$ws.Cells.Item(25, 9).Value = 1.032682321595814
$ws.Cells.Item(25, 10).Value = 1.042128583874373
$ws.Cells.Item(25, 11).Value = 1.040985968440215
$ws.Cells.Item(25, 12).Value = 1.0527813541608
$ws.Cells.Item(25, 13).Value = 1.05942378764956
$ws.Cells.Item(25, 14).Value = 1.043608525997294

